$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16:M18").Insert()
$ws.Range("F16:F18").Clear()

$ws.Range("B16").Value = "Driver's License ID"
$ws.Range("C16").Value = "Driver License ID"
$ws.Range("E16").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/j:PersonAugmentation/j:DriverLicense/j:DriverLicenseCardIdentification/nc:IdentificationID"

$ws.Range("B17").Value = "Driver License Source"
$ws.Range("E17").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/j:PersonAugmentation/j:DriverLicense/j:DriverLicenseCardIdentification/nc:IdentificationSourceText"

$ws.Range("B18").Value = "FBI ID"
$ws.Range("C18").Value = "FBI ID"
$ws.Range("E18").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/j:PersonAugmentation/j:PersonFBIIdentification/nc:IdentificationID"

$ws.Range("B17").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = "Driver License Source"

$ws.Rows("16:18").RowHeight = 56
